$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.527.07'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.901.79'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.21%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.26'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.98%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.909.11'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.97'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.360'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.409.00'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.536.29'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.68'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.906.52'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.61%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.99'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.14'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.56'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.72'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.57'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.08%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.84'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0843'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -9.84%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.68'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.62'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.17'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.68%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.59'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -6.31%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.76'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.48'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.72'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.290.45'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.648'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0583'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.48'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.94%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.96'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0238'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0922'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '248.85'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.30%  '
